# Apply updated crypto price / 1h-volume-change values to Sheet1.
# Values are plain text (not numbers), so NumberFormat is forced to
# "@" (Text) before the assignment to stop Excel from re-typing strings
# like "1.00" / "209.79" as numeric, then the style is reset back to
# "Normal" so no stray formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '28.354.25'
    'E2' = '  -0.90%  '
    'D3' = '1.549.81'
    'E3' = '  -2.02%  '
    'E4' = '  -0.05%  '
    'D5' = '209.79'
    'E5' = '  -1.77%  '
    'E6' = '  -1.61%  '
    'D7' = '1.00'
    'E7' = '  -0.09%  '
    'D8' = '23.88'
    'E8' = '  -0.97%  '
    'D10' = '0.0583'
    'E10' = '  -1.53%  '
    'D11' = '0.0889'
    'E11' = '  -0.40%  '
    'D12' = '1.770.42'
    'E12' = '  -2.02%  '
    'E13' = '  -1.66%  '
    'D14' = '28.328.63'
    'E14' = '  -1.02%  '
    'E15' = '  -2.32%  '
    'D16' = '0.510'
    'E16' = '  -1.92%  '
    'D17' = '60.89'
    'E17' = '  -2.13%  '
    'D18' = '227.80'
    'E18' = '  -1.66%  '
    'D19' = '7.35'
    'E19' = '  -0.72%  '
    'E20' = '  -2.26%  '
    'E21' = '  -0.01%  '
    'E22' = '  -0.14%  '
    'D23' = '8.92'
    'E23' = '  -2.90%  '
    'D24' = '2.03'
    'E24' = '  -2.45%  '
    'D25' = '151.39'
    'E25' = '  +0.19%  '
    'D26' = '14.74'
    'E27' = '  -0.99%  '
    'E28' = '  -0.09%  '
    'D29' = '6.24'
    'E29' = '  -3.23%  '
    'E30' = '  -3.45%  '
    'E31' = '  -4.45%  '
    'E32' = '  -1.34%  '
    'D33' = '1.383.87'
    'E33' = '  -0.93%  '
    'E34' = '  -3.45%  '
    'E35' = '  +0.94%  '
    'E36' = '  -3.59%  '
    'D37' = '2.34'
    'E37' = '  -1.13%  '
    'E38' = '  -3.00%  '
    'D39' = '0.0162'
    'E39' = '  -2.58%  '
    'D40' = '1.93'
    'E40' = '  +2.22%  '
    'D41' = '0.509'
    'E41' = '  -2.43%  '
    'E42' = '  -0.08%  '
    'D43' = '0.774'
    'E43' = '  -2.47%  '
    'D44' = '0.0455'
    'E44' = '  -2.67%  '
    'D45' = '5.37'
    'E45' = '  -1.64%  '
    'D46' = '61.98'
    'E46' = '  -2.11%  '
    'D47' = '1.683.91'
    'E47' = '  -2.06%  '
    'D48' = '0.870'
    'E48' = '  -9.49%  '
    'D49' = '85.75'
    'E49' = '  -0.95%  '
    'D50' = '42.20'
    'E50' = '  +4.56%  '
    'D51' = '0.0₆0103'
    'E51' = '  +0.43%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    $cell.Style = "Normal"
}
